$wb = $excel.ActiveWorkbook
$xlPasteFormats = [Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats

# --- Update the generation Date on the Metadata sheet ---
$meta = $wb.Worksheets.Item("Metadata")
$meta.Cells.Item(8, 2).Value2 = "2025-07-29T07:08:53+00:00"

# --- Elements sheet: widen a few columns to fit the new content ---
$ws = $wb.Worksheets.Item("Elements")
$ws.Columns.Item(1).ColumnWidth = 35.8828125
$ws.Columns.Item(2).ColumnWidth = 35.8828125
$ws.Columns.Item(11).ColumnWidth = 75.46875

# --- Insert a new row 6 for "CompetenceExclusive.exerciceProfessionnel", pushing
#     the existing "CompetenceExclusive.competenceExclusive" row down to row 7 ---
$ws.Range("A6:AJ6").Insert()

# Row 7 now holds exactly what row 6 used to hold (values + formatting). Re-apply
# that same formatting (fill/border/alignment) to the freshly inserted row 6.
$ws.Range("A7:AJ7").Copy()
$ws.Range("A6:AJ6").PasteSpecial($xlPasteFormats)

# --- Populate row 6 with the new "exerciceProfessionnel" element ---
$ws.Cells.Item(6, 1).Value2 = "CompetenceExclusive.exerciceProfessionnel"
$ws.Cells.Item(6, 2).Value2 = "CompetenceExclusive.exerciceProfessionnel"
$ws.Cells.Item(6, 11).Value2 = "Reference(https://interop.esante.gouv.fr/ig/fhir/mos/StructureDefinition/ExerciceProfessionnel)`n"
$ws.Cells.Item(6, 12).Value2 = "Lien vers la classe ExerciceProfessionnel."
$ws.Cells.Item(6, 13).Value2 = "Lien vers la classe ExerciceProfessionnel."
$ws.Cells.Item(6, 32).Value2 = "SavoirFaire.exerciceProfessionnel"

# Columns whose target text is empty ("") or looks purely numeric ("1") would be
# auto-typed as blank/number by Excel; force text via a quote-prefixed assignment,
# then strip the resulting quote-prefix style back to the plain text style.
$ws.Cells.Item(6, 4).Value2 = "'"
$ws.Cells.Item(6, 6).Value2 = "'1"
$ws.Cells.Item(6, 7).Value2 = "'1"
$ws.Cells.Item(6, 8).Value2 = "'"
$ws.Cells.Item(6, 9).Value2 = "'"
$ws.Cells.Item(6, 10).Value2 = "'"
$ws.Cells.Item(6, 16).Value2 = "'"
$ws.Cells.Item(6, 18).Value2 = "'"
$ws.Cells.Item(6, 19).Value2 = "'"
$ws.Cells.Item(6, 20).Value2 = "'"
$ws.Cells.Item(6, 21).Value2 = "'"
$ws.Cells.Item(6, 22).Value2 = "'"
$ws.Cells.Item(6, 23).Value2 = "'"
$ws.Cells.Item(6, 24).Value2 = "'"
$ws.Cells.Item(6, 25).Value2 = "'"
$ws.Cells.Item(6, 26).Value2 = "'"
$ws.Cells.Item(6, 27).Value2 = "'"
$ws.Cells.Item(6, 28).Value2 = "'"
$ws.Cells.Item(6, 29).Value2 = "'"
$ws.Cells.Item(6, 30).Value2 = "'"
$ws.Cells.Item(6, 31).Value2 = "'"
$ws.Cells.Item(6, 33).Value2 = "'1"
$ws.Cells.Item(6, 34).Value2 = "'1"
$ws.Cells.Item(6, 35).Value2 = "'"
$ws.Cells.Item(6, 36).Value2 = "'"

$fmtSrc = $ws.Range("D2")
$quotedCols = @(4, 6, 7, 8, 9, 10, 16, 18, 19, 20, 21, 22, 23, 24, 25, 26, 27, 28, 29, 30, 31, 33, 34, 35, 36)
foreach ($col in $quotedCols) {
    $fmtSrc.Copy()
    $ws.Cells.Item(6, $col).PasteSpecial($xlPasteFormats)
}
